$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "에피소드#06 - Pandas 데이터 전처리, 추가, 삭제, 데이터 type 변환"
$ws.Range("E4").Value = "https://teddylee777.github.io/pandas/pandas-tutorial-06"

$ws.Range("D24").Value = "[논문 요약 2021-02] RocketQA: An Optimized Training Approach to Dense Passage Retrieval"
$ws.Range("E24").Value = "https://blog.naver.com/hist0134/222215555386"

$ws.Range("D28").Value = "[keep9oing] Counterfactual Multi-Agent Policy Gradient (COMA) 리뷰 - (미완)"
$ws.Range("E28").Value = "https://ropiens.tistory.com/74"

$ws.Range("D36").Value = "Towards Contrastive Learning"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/308"

$ws.Range("D40").Value = "error_occured"
$ws.Range("E40").Value = "https://www.error_link.com"
